$d = $word.ActiveDocument

# Locate the run containing "proiect)" -> word changes to "proces)".
# "proiect)" = "pro" (unchanged) + "iect" (-> "ces") + ")" (unchanged text,
# but must become its own run, same as the target revision).
$found = $d.Content
$found.Find.Execute("proiect)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$wholeStart = $found.Start   # start of "proiect)"
$wholeEnd   = $found.End     # end of "proiect)"
$origSize  = $found.Font.Size # original (real) point size of this run

$prefixLen = 3                # length of "pro"
$midLen    = 4                # length of "iect" (being replaced by "ces")

$midStart   = $wholeStart + $prefixLen
$midEnd     = $midStart + $midLen
$parenStart = $midEnd
$parenEnd   = $wholeEnd

# Give the trailing ")" a distinct temporary format so it does not get
# silently re-absorbed into the run being edited just to the left of it.
$parenRange = $d.Range($parenStart, $parenEnd)
$parenRange.Font.Size = $origSize + 14

# Replace "iect" with "ces", tagging it with a distinguishing temporary
# format so the edit lands in its own run instead of merging back into
# the untouched "pro" run that precedes it.
$midRange = $d.Range($midStart, $midEnd)
$midRange.Font.Size = $origSize + 16
$midRange.Text = "ces"

# "ces" now occupies ($midStart, $midStart + 3); restore its real size.
$newMidEnd = $midStart + 3
$midFix = $d.Range($midStart, $newMidEnd)
$midFix.Font.Size = $origSize

# The ")" now starts right after "ces".
$parenStart2 = $newMidEnd
$parenEnd2 = $parenStart2 + 1

# Force the ")" into its own freshly-edited run (so it no longer carries
# the original run's identity) by rewriting its text via a throwaway
# placeholder character, then restoring it.
$parenRange2 = $d.Range($parenStart2, $parenEnd2)
$parenRange2.Font.Size = $origSize + 20
$parenRange2.Text = "Z"

$parenRange3 = $d.Range($parenStart2, $parenEnd2)
$parenRange3.Text = ")"

$parenFix = $d.Range($parenStart2, $parenEnd2)
$parenFix.Font.Size = $origSize
